# Update gh-pages output data (generated at 456a3b4)
$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 1036
$ws1.Range("G6").Value = 139
$ws1.Range("F8").Value = 198
$ws1.Range("F9").Value = 375
$ws1.Range("F14").Value = 12329
$ws1.Range("F15").Value = 85
$ws1.Range("F16").Value = 5472

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 1036
$ws4.Range("G8").Value = 139
$ws4.Range("F10").Value = 198
$ws4.Range("F16").Value = 12329
$ws4.Range("F18").Value = 85
$ws4.Range("F19").Value = 5472
